# Update calibration values across the risk-driver-type PnL factor sheets
# (AR, SETAR, GARCH, TARCH, AR_TARCH) to reflect the finalized virtual
# environment definition.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.000938174311667381
$ws.Range("B3").Value = 0.7809438713439286
$ws.Range("B4").Value = 0.1183709823931952
$ws.Range("B5").Value = "[0.9999999999999999, 0.2848119991792331, 0.2334356714793047, 0.1884336348211745, 0.23232495920002785, 0.4084685571727982, 0.24390892954114599, 0.24637731116838943, 0.2195950669032209, 0.19049379276934764, 0.18844228261075033, 0.168744803484475, 0.2220131456710551, 0.20815171860752224, 0.2130526454436412, 0.20095992871290813, 0.17677001340986137, 0.15930317633191526, 0.1790359238292295, 0.18302755053862838]"

# ---------------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value = 0.8918518518518519
$ws.Range("B4").Value = -0.1201920800837948
$ws.Range("B5").Value = 0.7146600096275811
$ws.Range("B6").Value = 0.08437256151599415
$ws.Range("B7").Value = 0.135949515535049
$ws.Range("B8").Value = 0.6457680558412187
$ws.Range("B9").Value = 0.06696972555906473
$ws.Range("B10").Value = "[1.0, 0.2709540677350821, 0.26721990309133725, 0.23072343204292584, 0.22874811051070423, 0.27251636979825433, 0.20516276535812522, 0.19685093437325013, 0.18630407918150968, 0.17245986290762602, 0.19478444685645047, 0.16961726962247356, 0.2011101630896298, 0.17966326359570117, 0.1823651825479308, 0.17906681112390138, 0.1798760075404138, 0.17708847267990474, 0.18537518343687875, 0.16282967064087508]"

# ---------------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.0004208491891122551
$ws.Range("B3").Value = 0.0003309222042712687
$ws.Range("B4").Value = 0.05181414205162736
$ws.Range("B5").Value = 0.9470128795253214
$ws.Range("B6").Value = "[1.0, 0.012921484748652453, -0.026072343461510026, -0.06395318896691948, -0.07457061134478891, 0.28415198586866247, 0.011226160866768017, 0.042061986071826055, 0.0074122145730563015, -0.035960630026002874, -0.03897814043332797, -0.058845018297552946, -0.008157874727687857, -0.008459946571520601, 0.022237295878262855, 0.019433928941365638, -0.017960817765663915, -0.04639642287492738, -0.07529952439314666, -0.012642632212823347]"

# ---------------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.001562538384524367
$ws.Range("B3").Value = 0.0003317251781729857
$ws.Range("B4").Value = 0.04650875848957307
$ws.Range("B5").Value = 0.946142285572969
$ws.Range("B6").Value = "[1.0, 0.01170640069309429, -0.027019871792121383, -0.06473556769136818, -0.07523460134067694, 0.28394915683882177, 0.010571844304318033, 0.04139500528711297, 0.006843809719162673, -0.036535031037552776, -0.03956013651383042, -0.05902292178768356, -0.008589463593761677, -0.008996289543797251, 0.02123524169535309, 0.018803478144489818, -0.01825499053900965, -0.04635009855593515, -0.07558730364559292, -0.012818501239923279]"
$ws.Range("B7").Value = 0.0125454530356663

# ---------------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.00319136883951554
$ws.Range("B3").Value = 0.0002683582271512838
$ws.Range("B4").Value = 0.04218000568309339
$ws.Range("B5").Value = 0.9502592922823111
$ws.Range("B6").Value = "[0.9999999999999998, 0.03436747043835898, -0.009916977392808952, -0.036342088441116874, -0.03398037899371765, 0.2278832559768839, -0.015526867512406386, 0.024832353317328033, -0.0012328249086491533, -0.04742429740295087, -0.021595739425330534, -0.049224940512919396, -0.0010460084428602415, -0.006020314176895672, 0.007707317433711263, 0.02113771040858326, -0.019646713496766777, -0.04458143120438001, -0.04379339936949635, -0.014302923625889639]"
$ws.Range("B7").Value = 0.01257965920255122
$ws.Range("B9").Value = 0.7871140150637183
